$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Objetivos:" value changes to the docente name
$ws.Range("B10").Value = '198273 - Domingos Savio Giordani'
$ws.Range("C10").Value = '198273 - Domingos Savio Giordani'

# Row 13: add label cell A13 (copy style from A14), update B13/C13 text, set row height 60
$ws.Range("A14").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

# Row 14: label + long text change (row height stays 60)
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Training and team work, Communication, Systematic Innovation, Legislation, Project Management. Problem Identification, Project Formulation, Problem Specification, Available Knowledge Analysis, Evaluation and Decision Making, Timeline, Reporting, Project Submission'
$ws.Range("C14").Value = 'Training and team work, Communication, Systematic Innovation, Legislation, Project Management. Problem Identification, Project Formulation, Problem Specification, Available Knowledge Analysis, Evaluation and Decision Making, Timeline, Reporting, Project Submission'

# Row 15: label change; B15/C15 reuse the "01/01/2020" text (copy as value from B8/C8 to avoid date auto-parsing); row height 60 -> 120
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Rows.Item(15).RowHeight = 120

# Row 16: label + long text change (row height stays 120)
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = 'Training and work in teams and communication - the development of skills essential to work in teams; Systematic Innovation - development of innovative solutions, systematization and characteristics; Legislation - notions of legislation applied to corporate action; Project Management and Schedule - Methodologies and necessary schematizations with the management elements; Problem Identification - systematization of actions to locate causes; Formulation of the Project - presentation of the managerial aspects necessary for the development of the project, Management Plan, Project Analytical Structure (EAP) etc; Specification of Problems - systematization of problems within the areas of knowledge; Analysis of Available Knowledge, Evaluation and Decision Making; Reporting - formatting within ABNT standards; Presentation of Projects.'
$ws.Range("C16").Value = 'Training and work in teams and communication - the development of skills essential to work in teams; Systematic Innovation - development of innovative solutions, systematization and characteristics; Legislation - notions of legislation applied to corporate action; Project Management and Schedule - Methodologies and necessary schematizations with the management elements; Problem Identification - systematization of actions to locate causes; Formulation of the Project - presentation of the managerial aspects necessary for the development of the project, Management Plan, Project Analytical Structure (EAP) etc; Specification of Problems - systematization of problems within the areas of knowledge; Analysis of Available Knowledge, Evaluation and Decision Making; Reporting - formatting within ABNT standards; Presentation of Projects.'

# Row 17: label changes, clear B17/C17 entirely, reset row height to default
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17:C17").Clear()
$ws.Rows.Item(17).AutoFit()

# Row 18: label changes, add B18/C18 (copy style from B19/C19), set row height 60
$ws.Range("A18").Value = 'Método:'
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("B18").Value = '198273 - Domingos Savio Giordani'
$ws.Range("C18").Value = '198273 - Domingos Savio Giordani'
$ws.Rows.Item(18).RowHeight = 60

# Row 19: label change only, text stays the same, row height stays 60
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Apresentações intermediárias e finais.'
$ws.Range("C19").Value = 'Apresentações intermediárias e finais.'

# Row 20: label change only, text stays the same, row height stays 60
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Serão feitas duas avaliações por uma banca de professores que assistirão às apresentações, as notas serão as médias das notas dadas pelos professores.'
$ws.Range("C20").Value = 'Serão feitas duas avaliações por uma banca de professores que assistirão às apresentações, as notas serão as médias das notas dadas pelos professores.'

# Row 21: label change, text stays the same, row height 60 -> 120
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina.'
$ws.Range("C21").Value = 'Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina.'
$ws.Rows.Item(21).RowHeight = 120

# Row 22: no longer needed, delete it entirely (shrinks dimension to C21)
$ws.Rows.Item(22).Delete()
